$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the two hyperlink paragraphs (themuse.com / collegeinfogeek.com)
#    plus the blank paragraph that separates them. We locate them via the
#    document's Hyperlinks collection, which is more robust than fixed
#    paragraph indices.
# ------------------------------------------------------------------
$muse = $null
$geek = $null
foreach ($h in $d.Hyperlinks) {
    if ($h.Address -like "*themuse.com*") { $muse = $h }
    if ($h.Address -like "*collegeinfogeek.com*") { $geek = $h }
}

if (($muse -ne $null) -and ($geek -ne $null)) {
    $startRange = $d.Range($muse.Range.Start, $muse.Range.End)
    $startRange.Expand(4) | Out-Null
    $endRange = $d.Range($geek.Range.Start, $geek.Range.End)
    $endRange.Expand(4) | Out-Null

    $killRange = $d.Range($startRange.Start, $endRange.End)
    $killRange.Delete()
}

# ------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark: it currently sits alone in its
#    own (empty) paragraph right after "Contact - email, LinkedIn,
#    Instagram"; it should instead mark the very start of the
#    "- Navigation" paragraph. Deleting that now-empty paragraph's
#    range removes both the paragraph and the bookmark it contains.
#    (Guarded so this only fires when the bookmark still lives alone in
#    its own blank paragraph, matching the document's pre-edit shape.)
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBmRange = $d.Range($oldBm.Start, $oldBm.End)
    $oldBmRange.Expand(4) | Out-Null
    if ($oldBmRange.Text -eq "`r") {
        $oldBmRange.Delete()
    }
}

# ------------------------------------------------------------------
# 3) Re-create the "_GoBack" bookmark, collapsed, at the start of the
#    "- Navigation" paragraph.
# ------------------------------------------------------------------
$navRange = $d.Content
$navRange.Find.Execute("- Navigation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$navStart = $navRange.Start
$bmRange = $d.Range($navStart, $navStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
